$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths for newly added columns L:P (18 chars, matching existing F:K) ---
$ws.Range("L1:P1").ColumnWidth = 17.17

# --- Unmerge ranges that need to grow before we touch their cells ---
$ws.Range("F1:K1").UnMerge()
$ws.Range("K2").UnMerge()

# --- Row 2: month headers ---
$ws.Range("F2").Value = "February"
$ws.Range("K2").Value = "March"
$ws.Range("O2").Value = "April"

# --- Row 3: weekly date ranges (Feb - Apr), style s3 (was s4) ---
$ws.Range("F3").Value = "01/Feb - 07/Feb"
$ws.Range("G3").Value = "08/Feb - 14/Feb"
$ws.Range("H3").Value = "15/Feb - 21/Feb"
$ws.Range("I3").Value = "22/Feb - 28/Feb"
$ws.Range("J3").Value = "29/Feb - 06/Mar"
$ws.Range("K3").Value = "07/Mar - 13/Mar"
$ws.Range("L3").Value = "14/Mar - 20/Mar"
$ws.Range("M3").Value = "21/Mar - 27/Mar"
$ws.Range("N3").Value = "28/Mar - 03/Apr"
$ws.Range("O3").Value = "04/Apr - 10/Apr"
$ws.Range("P3").Value = "11/Apr - 17/Apr"

# new week-header cells need the same wrap/center style as F3:K3 (style index "3" post-edit == fontId2/fillId2/center)
$ws.Range("F3:P3").HorizontalAlignment = -4108
$ws.Range("F3:P3").Font.Color = $ws.Range("F2").Font.Color
$ws.Range("F3:P3").Interior.Color = $ws.Range("F2").Interior.Color

# --- Re-merge the grown month/year headers ---
$ws.Range("F1:P1").Merge()
$ws.Range("K2:N2").Merge()
$ws.Range("O2:P2").Merge()

# --- Clear Activity names (column C) for existing rows; keep Tasks numbering/Start/End dates ---
$ws.Range("C4:C8").ClearContents()

# --- Update existing rows 4-8 (Start/End dates + move the colored marker cell) ---
$ws.Range("D4").Value = "02/01"
$ws.Range("E4").Value = "02/07"

$ws.Range("D5").Value = "02/08"
$ws.Range("E5").Value = "02/14"

$ws.Range("D6").Value = "02/15"
$ws.Range("E6").Value = "02/21"

$ws.Range("D7").Value = "02/22"
$ws.Range("E7").Value = "02/28"

$ws.Range("D8").Value = "02/29"
$ws.Range("E8").Value = "03/06"

# --- Add new rows 9-13, copying formatting down from row 8 ---
$ws.Range("B8:E8").Copy()
$ws.Range("B9:E13").PasteSpecial(-4122)
$ws.Range("F4").Copy()
$ws.Range("K9").PasteSpecial(-4122)
$ws.Range("L10").PasteSpecial(-4122)
$ws.Range("M11").PasteSpecial(-4122)
$ws.Range("N12").PasteSpecial(-4122)
$ws.Range("O13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B9").Value = 6
$ws.Range("C9").ClearContents()
$ws.Range("D9").Value = "03/07"
$ws.Range("E9").Value = "03/13"

$ws.Range("B10").Value = 7
$ws.Range("C10").ClearContents()
$ws.Range("D10").Value = "03/14"
$ws.Range("E10").Value = "03/20"

$ws.Range("B11").Value = 8
$ws.Range("C11").ClearContents()
$ws.Range("D11").Value = "03/21"
$ws.Range("E11").Value = "03/27"

$ws.Range("B12").Value = 9
$ws.Range("C12").ClearContents()
$ws.Range("D12").Value = "03/28"
$ws.Range("E12").Value = "04/03"

$ws.Range("B13").Value = 10
$ws.Range("C13").ClearContents()
$ws.Range("D13").Value = "04/04"
$ws.Range("E13").Value = "04/10"
